$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Fitness" column (C) values for the ranges of rows whose
# values changed, per the commit diff.
$ws.Range("C2:C10").Value = 10206
$ws.Range("C11:C13").Value = 8622
$ws.Range("C14:C15").Value = 8521
$ws.Range("C16:C18").Value = 8345
$ws.Range("C19:C21").Value = 8244
$ws.Range("C22:C28").Value = 8055
$ws.Range("C100:C252").Value = 7293
